$d = $word.ActiveDocument

# 1. Rewrite the opening paragraph ("This week, we will be taking a look...")
#    into the new continuation-themed intro paragraph. Using Find/Replace so
#    that the two existing runs get collapsed into the single new run.
$old1 = "This week, we will be taking a look at how we can create a tea cup using the Curve Pen tool. You can only have access to this tool in edit mode, if you have added a curved object in Object mode first. If you have something like a cube, you will never see this tool"
$new1 = "This week will be a continuation of the last tutorial that we did. Last week, we created the right side of a tea cup using the Curve Pen tool. This week we will be forming the handle that will be attached to the cup. And again, we will be using the Curve Pen Tool to achieve this."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2. Rewrite the "So, if you are interested in..." paragraph into the new
#    "So, if this sounds at all interesting..." paragraph (also collapses runs).
$old2 = "So, if you are interested in learning an alternative method to churning out some curves. Then please join us for our brand-new article this week entitled:"
$new2 = "So, if this sounds at all interesting to you than please join us for our brand-new article entitled:"
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3. Insert a brand-new Heading1 paragraph right after that paragraph, holding
#    the new article title "The Curve Pen Tool Part 2".
$introPara = $d.Paragraphs(3)
$introPara.Range.InsertParagraphAfter()
$titlePara = $d.Paragraphs(4)
$titlePara.Range.Text = "The Curve Pen Tool Part 2"
$titlePara.Style = "Heading1"

# 4. Remove the old Heading1 paragraph further down ("9 The Curve Pen Tool").
$old3 = "9 The Curve Pen Tool"
$searchRange = $d.Content
$searchRange.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($searchRange.Find.Found) {
    $oldHeadingPara = $searchRange.Paragraphs(1)
    $oldHeadingPara.Range.Delete()
}
